$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 68-76: these rows had their value in column F ("FVTPL value" area spilled
# into F because no "Quarter" ticker was recorded). The upload adds a "PBT"
# quarter tag in column C and moves the numeric value from F into D.
$rows = 68..76

foreach ($r in $rows) {
    $fCell = $ws.Cells.Item($r, 6)   # column F
    $val = $fCell.Value2
    $dCell = $ws.Cells.Item($r, 4)   # column D
    $ws.Cells.Item($r, 3).Value2 = "PBT"   # column C
    $dCell.Value2 = $val
    $dCell.NumberFormat = "#,##0"   # same numeric style the F cells used
    $fCell.Clear()                  # fully remove the now-unused F cell
}

# Update the sheet view: selection moves to D1 and the previous scroll
# position (topLeftCell = A58) is cleared.
$ws.Range("D1").Select()
